# "aggiornamento fino a 28 luglio" - append daily rows 302:328 (dates 44376..44402)
# to the Soliera sheet, continuing the existing A:D table (date, nuovi pos.,
# somma mobile 7gg., somma mobile 7gg. per 100mila abitanti).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @(
    # date(serial), nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.
    @(44376, 1, 1, 6.458696635019054),
    @(44377, 0, 1, 6.458696635019054),
    @(44378, 0, 1, 6.458696635019054),
    @(44379, 0, 1, 6.458696635019054),
    @(44380, 0, 1, 6.458696635019054),
    @(44381, 0, 1, 6.458696635019054),
    @(44382, 0, 1, 6.458696635019054),
    @(44383, 0, 0, 0),
    @(44384, 0, 0, 0),
    @(44385, 0, 0, 0),
    @(44386, 0, 0, 0),
    @(44387, 0, 0, 0),
    @(44388, 0, 0, 0),
    @(44389, 0, 0, 0),
    @(44390, 0, 0, 0),
    @(44391, 0, 0, 0),
    @(44392, 0, 0, 0),
    @(44393, 0, 0, 0),
    @(44394, 0, 0, 0),
    @(44395, 0, 0, 0),
    @(44396, 1, 1, 6.458696635019054),
    @(44397, 0, 1, 6.458696635019054),
    @(44398, 0, 1, 6.458696635019054),
    @(44399, 1, 2, 12.91739327003811),
    @(44400, 1, 3, 19.37608990505716),
    @(44401, 1, 4, 25.83478654007622),
    @(44402, 0, 4, 25.83478654007622)
)

$startRow = 302
for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $startRow + $i
    $row = $newData[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
$endRow = $startRow + $newData.Count - 1

# Column A carries the bold/centered/bordered date style (same as the rest of
# the column) - copy formats only from the last existing dated row so we
# reuse the existing style record instead of minting a new one.
$ws.Range("A301").Copy()
$ws.Range("A" + $startRow + ":A" + $endRow).PasteSpecial(-4122)
